# Generate Report for Archive
#
# 1. Shared string "Ready for handoff" -> "In Translation" everywhere it is
#    used (Overview!E2:F2, zh-cn!C2, de-de!C2 — all four cells hold the same
#    status text for the current localization stage).
# 2. Narrow the "handoff status" columns (Overview E:F, and column C on the
#    per-locale sheets) from ~17.22 chars to ~13.41 chars now that the
#    shorter "In Translation" label no longer needs the extra room.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Range("E1:F1").ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Range("C1").ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Range("C1").ColumnWidth = 12.5
